# TianeTi.xlsx edit:
#  1. Rename the only worksheet from "1" to "Tianeti".
#  2. Delete the blank row 8 (the "Note" row then shifts up from row 9 to row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet "1" -> "Tianeti"
$ws.Name = "Tianeti"

# Remove the empty row 8, shifting subsequent rows (the Note row) up by one.
$ws.Rows("8").Delete()
